$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Added known GUI defects to StatusSheet"
#
# The author added three newly-discovered GUI defects to the
# "Defects Found" sheet (Name / Description / Section ID / Category).
# Cell values are entered in the same order the author typed them so
# that the shared-string table comes out in the same sequence: the
# repeated "GUI" category first, then the Logout row, then Color, then
# Validation.
# ---------------------------------------------------------------------

$ws = $wb.Worksheets.Item("Defects Found")

$ws.Range("D2").Value = "GUI"
$ws.Range("D3").Value = "GUI"
$ws.Range("D4").Value = "GUI"

$ws.Range("A4").Value = "Logout Functionality"
$ws.Range("B4").Value = """Log Out"" button on GUI does not perform any functions."

$ws.Range("A2").Value = "Color Requirement"
$ws.Range("B2").Value = "Requirement of blue and green GUI is not met. "

$ws.Range("A3").Value = "Validation of Additions"
$ws.Range("B3").Value = "When administrators add new users or locations, there is no notification of how duplicates are handled from the client side."

$ws.Range("C2").Value = 2.1
$ws.Range("C3").Value = 2.2
$ws.Range("C4").Value = 3.1

# Header formatting: "Section ID" / "Category" columns centered like
# the existing "Name" column; "Description" stays bold but switches to
# a left/top wrapping layout to fit the longer defect descriptions.
$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("A1").VerticalAlignment = -4108
$ws.Range("C1").HorizontalAlignment = -4108
$ws.Range("C1").VerticalAlignment = -4108
$ws.Range("D1").HorizontalAlignment = -4108
$ws.Range("D1").VerticalAlignment = -4108

$ws.Range("B1").HorizontalAlignment = -4131
$ws.Range("B1").VerticalAlignment = -4160
$ws.Range("B1").WrapText = $true

# Data rows: Name / Section ID / Category centered; Description
# left-aligned, top-aligned and wrapped so long text is readable.
$ws.Range("A2:A4").HorizontalAlignment = -4108
$ws.Range("A2:A4").VerticalAlignment = -4108
$ws.Range("C2:C4").HorizontalAlignment = -4108
$ws.Range("C2:C4").VerticalAlignment = -4108
$ws.Range("D2:D4").HorizontalAlignment = -4108
$ws.Range("D2:D4").VerticalAlignment = -4108

$ws.Range("B2:B4").HorizontalAlignment = -4131
$ws.Range("B2:B4").VerticalAlignment = -4160
$ws.Range("B2:B4").WrapText = $true

# Row heights sized to fit the wrapped description text.
$ws.Rows.Item(2).RowHeight = 48.75
$ws.Rows.Item(3).RowHeight = 66.75
$ws.Rows.Item(4).RowHeight = 31.5

# Column widths sized for the new "Name"/"Description" content.
$ws.Columns.Item(1).ColumnWidth = 19.5
$ws.Columns.Item(2).ColumnWidth = 26.833333333333336
$ws.Columns.Item(3).ColumnWidth = 10.166666666666666
$ws.Columns.Item(4).ColumnWidth = 10.166666666666666

$ws.PageSetup.Orientation = 1

# The author had been working in "Defect Checklist" and switched to
# "Defects Found" after entering the new rows, leaving the cursor a
# couple of rows below the freshly-typed data.
$checklist = $wb.Worksheets.Item("Defect Checklist")
$checklist.Range("D8").Select()

$ws.Activate()
$ws.Range("D7").Select()
